$d = $word.ActiveDocument
$p8 = $d.Paragraphs(8)
$full = $p8.Range.Duplicate
$full.MoveEnd(1, -1)
$full.Text = "Distributive | Kingston, Ontario, Canada"
$p8b = $d.Paragraphs(8)
$base = $p8b.Range.Start

function DumpP8 {
  param($label)
  $d2 = $word.ActiveDocument
  $pp = $d2.Paragraphs(8)
  Write-Host "=== $label === text=[$($pp.Range.Text)]"
}

DumpP8 "after full text set"

$seg = $p8b.Range.Duplicate
$seg.SetRange($base+12, $base+13)
$seg.Font.Bold = $false
DumpP8 "after seg1 bold"

$seg2 = $p8b.Range.Duplicate
$seg2.SetRange($base+13, $base+14)
$seg2.Font.Bold = $false
DumpP8 "after seg2 bold"
$seg2.Font.Color = 12040119
DumpP8 "after seg2 color"

$seg3 = $p8b.Range.Duplicate
$seg3.SetRange($base+14, $base+15)
$seg3.Font.Bold = $false
DumpP8 "after seg3 bold"

$seg4 = $p8b.Range.Duplicate
$seg4.SetRange($base+15, $base+41)
$seg4.Font.Size = 10
DumpP8 "after seg4 size"
$seg4.Font.Color = 6710886
DumpP8 "after seg4 color"
